$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E10").Value = 270
$ws.Range("E11").Value = 195
$ws.Range("E12").Value = 284
$ws.Range("E14").Value = 81
$ws.Range("E15").Value = 109
$ws.Range("E22").Value = 107

$ws.Range("E25").Value = 139
$ws.Range("F25").Value = 69
$ws.Range("H25").Value = 69

$ws.Range("E29").Value = 117
$ws.Range("F29").Value = 73
$ws.Range("H29").Value = 73

$ws.Range("E30").Value = 135

$ws.Range("E34").Value = 131
$ws.Range("F34").Value = 86
$ws.Range("H34").Value = 86

$ws.Range("E42").Value = 219
$ws.Range("E46").Value = 170
$ws.Range("E50").Value = 129
